# "week 9 all subjects added"
#
# Sheet1 is the weekly tracker: column C = TPL, column D = SRE, column E = ALGO.
# This commit:
#   - adds a "Slides:" line to the existing TPL Week#4 (row6/C6), SRE Week#8
#     (row11/D11) and TPL Week#8 (row11/C11) entries
#   - adds brand-new "Week 9" rows for every subject (row12), plus a TPL
#     "Week 10" entry (row13/E13)
#   - appends extra "Assignment ... Solved" / "Topics" lines to the existing
#     ALGO Week#8 (row11/E11) and ALGO Week#9 (row12/E12) entries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Popping directional-formatting mark used in the original "Kamran (Abasyn)"
# note text (U+202C).
$pdf = [char]0x202C

# --- Row 6: TPL Week 4 — add a Slides line before the Topics line ---------
$ws.Cells.Item(6, 3).Value = (
    "Topic: name- Week 4, lectures- 1 Lecture, duration- 00:43;`n" +
    "Video: link- https://drive.google.com/file/d/1JZjwQWO7_aTq04qMJLL5pf6HXNrv_qyN/preview, name- TPL Week#4, duration- 00:43;`n" +
    "Slides: slide- Chapter 4.ppt;`n" +
    "Topics: Compiler working flow chart, Top down parsing alogrithum, Construction of Syntaxtical Structure using BNF;"
)

# --- Row 11: TPL Week 8 — add a Slides line -------------------------------
$ws.Cells.Item(11, 3).Value = (
    "Topic: name- Week 8, lectures- , duration- ;`n" +
    "Slides: slide- Chapter 6.ppt;`n" +
    "Assignment: name- Assignment no 3 (Deadline 15 Dec 2019), img- /TPL/Assignment No 3 TPL.png;`n" +
    "Note: heading- Note, text- Video lectures not uploaded on google drive yet. If you want them faster please call Kamran (Abasyn) at 0302 5003156$pdf.;"
)

# --- Row 11: SRE Week 8 — add a Slides line -------------------------------
$ws.Cells.Item(11, 4).Value = (
    "Topic: name- Week 8, lectures- 3 Lectures, duration- 01:22;`n" +
    "Video: link- https://drive.google.com/file/d/18sP4Wf0WZh_3gh0PCDZaPH0YT7HsDmhn/preview, name- SRE Week8 Part#1.mp4, duration- 00:01;`n" +
    "Video: link- https://drive.google.com/file/d/1cFUxdRC5VwlsjxILvMbMqO6SVWkXwZB1/preview, name- SRE Week8 Part#2.mp4, duration- 00:51;`n" +
    "Video: link- https://drive.google.com/file/d/128EIgEoYNkqrN_PCG30nwpn1MGaNhK-J/preview, name- SRE Week8 Part#3.mp4, duration- 00:30;`n" +
    "Assignment: name- Assignment no 3 (Deadline 21 Dec 2019), img- /SRE/Assignment 3.png;`n" +
    "Slides: slide- Lecture Slide_5.pptx;"
)

# --- Row 11: ALGO Week 8 — extra solved-assignment + more topics ---------
$ws.Cells.Item(11, 5).Value = (
    "Topic: name- Week 8, lectures- 2 Lectures, duration- 01:33;`n" +
    "Video: link- https://drive.google.com/file/d/1lpPNJAvs6WzQuJ6z0dhX96HCxnarShVs/preview, name- AD&AA Week # 8 Part1.mp4, duration- 00:54;`n" +
    "Video: link- https://drive.google.com/file/d/1A8aBmhTcErBz6hETG1YYLN-GMCLMZUSX/preview, name- AD&AA Week # 8 Part2.mp4, duration- 00:39;`n" +
    "Assignment: name- Assignment 2 (Deadline 22 Nov), img- /ALGO/Assignment 2.png*/ALGO/Assignment 2 DL.png;`n" +
    "Assignment: name- Assignment 2 Solved, img- /ALGO/Assignment20.jpg*/ALGO/Assignment21.jpg;`n" +
    "Slides: slide- lec4.pptx;`n" +
    "Topics: Redex Sort, Bubble Sort, Merge Sort, Counting Sort, Bucket Sort;"
)

# --- Row 12: TPL Week 9 (brand new) ---------------------------------------
$ws.Cells.Item(12, 3).WrapText = $true
$ws.Cells.Item(12, 3).Font.Bold = $true
$ws.Cells.Item(12, 3).HorizontalAlignment = -4131
$ws.Cells.Item(12, 3).VerticalAlignment = -4108
$ws.Cells.Item(12, 3).Value = (
    "Topic: name- Week 9, lectures- , duration- ;`n" +
    "Slides: slide- Chapter 6.ppt;`n" +
    "Note: heading- Note, text- Video lectures not uploaded on google drive yet. If you want them faster please call Kamran (Abasyn) at 0302 5003156$pdf.;"
)

# --- Row 12: SRE Week 9 (brand new) ---------------------------------------
$ws.Cells.Item(12, 4).WrapText = $true
$ws.Cells.Item(12, 4).Font.Bold = $true
$ws.Cells.Item(12, 4).HorizontalAlignment = -4131
$ws.Cells.Item(12, 4).VerticalAlignment = -4108
$ws.Cells.Item(12, 4).Value = (
    "Topic: name- Week 9, lectures- , duration- ;`n" +
    "Slides: slide- Lecture Slide_5.pptx;`n" +
    "Important: Quiz from Lecture 5 (above lecture) in next class. Please come prepared.;"
)

# --- Row 12: ALGO Week 9 — extra solved-assignment + more topics ---------
$ws.Cells.Item(12, 5).Value = (
    "Topic: name- Week 9, lectures- 2 Lectures, duration- 01:51;`n" +
    "Video: link- https://drive.google.com/file/d/1430C-n3l2kRoy2Qn76kBep_umfClXmkF/preview, name- AD&AA Week # 9 Part1.mp4, duration- 00:38;`n" +
    "Video: link- https://drive.google.com/file/d/1UBidBALGcA7KB8t11RDEwlTti64BilJf/preview, name- AD&AA Week # 9 Part2.mp4, duration- 01:13;`n" +
    "Assignment: name- Assignment 3 (Deadline 29 Nov), img- /ALGO/Assignment 3.png;`n" +
    "Assignment: name- Assignment 3 Solved, img- /ALGO/Assignment30.jpg*/ALGO/Assignment31.jpg;`n" +
    "Slides: slide- lec5.ppt;`n" +
    "Important: Quiz in next class from lecture 5 (above slides);`n" +
    "Topics: Pigeon hole sort, Important points for algorithms, Chapter 5, Asymptotatic analysis, Rate of growth, Wilson, Fibonacci, Investing on algorithms or computing?;"
)

# --- Row 13: TPL Week 10 (brand new) --------------------------------------
$ws.Cells.Item(13, 5).WrapText = $true
$ws.Cells.Item(13, 5).Font.Bold = $true
$ws.Cells.Item(13, 5).HorizontalAlignment = -4131
$ws.Cells.Item(13, 5).VerticalAlignment = -4108
$ws.Cells.Item(13, 5).Value = (
    "Topic: name- Week 10, lectures- , duration- ;`n" +
    "Slides: slide- Lec6.pptx;`n" +
    "Note: heading- Video lectures not uploaded yet;"
)

# --- Row heights, to fit the new (taller) wrapped text --------------------
$ws.Rows.Item(11).RowHeight = 306
$ws.Rows.Item(12).RowHeight = 323
$ws.Rows.Item(13).RowHeight = 51

# --- Selection / scroll position, matching the saved view -----------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("C12").Select()

Write-Host "week 9 all subjects added"
